$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "69.073.08"
$ws.Cells.Item(3, 4).Value = "3.936.04"
$ws.Cells.Item(3, 5).Value = "  +0.75%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.29%  "
$ws.Cells.Item(5, 4).Value = "487.46"
$ws.Cells.Item(5, 5).Value = "  +0.60%  "
$ws.Cells.Item(6, 4).Value = "147.81"
$ws.Cells.Item(6, 5).Value = "  +1.51%  "
$ws.Cells.Item(7, 4).Value = "0.623"
$ws.Cells.Item(7, 5).Value = "  -0.23%  "
$ws.Cells.Item(8, 4).Value = "0.997"
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 5).Value = "  +1.68%  "
$ws.Cells.Item(10, 4).Value = "0.177"
$ws.Cells.Item(10, 5).Value = "  +6.08%  "
$ws.Cells.Item(11, 4).Value = "0.0000347"
$ws.Cells.Item(11, 5).Value = "  -1.60%  "
$ws.Cells.Item(12, 4).Value = "43.14"
$ws.Cells.Item(12, 5).Value = "  +1.74%  "
$ws.Cells.Item(13, 4).Value = "10.51"
$ws.Cells.Item(13, 5).Value = "  -0.90%  "
$ws.Cells.Item(14, 4).Value = "4.554.79"
$ws.Cells.Item(14, 5).Value = "  +0.48%  "
$ws.Cells.Item(15, 4).Value = "3.932.36"
$ws.Cells.Item(15, 5).Value = "  +0.60%  "
$ws.Cells.Item(16, 4).Value = "14.30"
$ws.Cells.Item(16, 5).Value = "  -2.63%  "
$ws.Cells.Item(17, 5).Value = "  -0.68%  "
$ws.Cells.Item(18, 4).Value = "20.01"
$ws.Cells.Item(18, 5).Value = "  +1.13%  "
$ws.Cells.Item(19, 5).Value = "  +1.87%  "
$ws.Cells.Item(20, 4).Value = "69.099.88"
$ws.Cells.Item(20, 5).Value = "  +1.42%  "
$ws.Cells.Item(21, 4).Value = "438.57"
$ws.Cells.Item(21, 5).Value = "  -2.00%  "
$ws.Cells.Item(22, 5).Value = "  +4.75%  "
$ws.Cells.Item(23, 4).Value = "14.69"
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 4).Value = "89.57"
$ws.Cells.Item(24, 5).Value = "  +0.72%  "
$ws.Cells.Item(25, 4).Value = "12.04"
$ws.Cells.Item(25, 5).Value = "  +8.35%  "
$ws.Cells.Item(26, 4).Value = "3.72"
$ws.Cells.Item(26, 5).Value = "  +3.41%  "
$ws.Cells.Item(27, 4).Value = "11.12"
$ws.Cells.Item(27, 5).Value = "  -4.27%  "
$ws.Cells.Item(28, 4).Value = "37.24"
$ws.Cells.Item(29, 4).Value = "5.65"
$ws.Cells.Item(29, 5).Value = "  -3.19%  "
$ws.Cells.Item(30, 4).Value = "712.80"
$ws.Cells.Item(30, 5).Value = "  +3.30%  "
$ws.Cells.Item(31, 5).Value = "  +1.16%  "
$ws.Cells.Item(32, 4).Value = "13.42"
$ws.Cells.Item(32, 5).Value = "  +0.53%  "
$ws.Cells.Item(33, 5).Value = "  +1.19%  "
$ws.Cells.Item(34, 4).Value = "0.471"
$ws.Cells.Item(34, 5).Value = "  +29.32%  "
$ws.Cells.Item(35, 4).Value = "0.0₃0914"
$ws.Cells.Item(35, 5).Value = "  -0.96%  "
$ws.Cells.Item(36, 4).Value = "6.06"
$ws.Cells.Item(36, 5).Value = "  +5.95%  "
$ws.Cells.Item(37, 2).Value = "OKB"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(37, 4).Value = "61.02"
$ws.Cells.Item(37, 5).Value = "  +3.74%  "
$ws.Cells.Item(38, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(38, 4).Value = "40.94"
$ws.Cells.Item(38, 5).Value = "  -1.55%  "
$ws.Cells.Item(39, 5).Value = "  -0.48%  "
$ws.Cells.Item(40, 4).Value = "1.00"
$ws.Cells.Item(40, 5).Value = "  +0.16%  "
$ws.Cells.Item(41, 5).Value = "  +0.00%  "
$ws.Cells.Item(42, 4).Value = "2.95"
$ws.Cells.Item(42, 5).Value = "  +1.32%  "
$ws.Cells.Item(43, 5).Value = "  +2.41%  "
$ws.Cells.Item(44, 4).Value = "3.09"
$ws.Cells.Item(44, 5).Value = "  +1.55%  "
$ws.Cells.Item(45, 5).Value = "  -0.02%  "
$ws.Cells.Item(46, 4).Value = "0.0₆0378"
$ws.Cells.Item(46, 5).Value = "  +15.52%  "
$ws.Cells.Item(47, 5).Value = "  +9.21%  "
$ws.Cells.Item(48, 5).Value = "  +0.93%  "
$ws.Cells.Item(49, 2).Value = "Stacks"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(49, 4).Value = "2.99"
$ws.Cells.Item(49, 5).Value = "  +6.15%  "
$ws.Cells.Item(50, 2).Value = "LidoDAOToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(50, 4).Value = "3.37"
$ws.Cells.Item(50, 5).Value = "  -1.12%  "
$ws.Cells.Item(51, 4).Value = "2.09"
$ws.Cells.Item(51, 5).Value = "  -1.61%  "
